$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - B2 and D2 updates
$ws.Range("B2").Value = 'Running Robot tests (Assigning) studying Assigning generously Assign all operators Assign all operators for the remaining users'
$ws.Range("D2").Value = '[''Running'', ''Robot'', ''tests'', ''('', ''Assigning'', '')'', ''studying'', ''Assigning'', ''generously'', ''Assign'', ''all'', ''operators'', ''Assign'', ''all'', ''operators'', ''for'', ''the'', ''remaining'', ''users'']'

# Row 15 - B15 and D15 updates
$ws.Range("B15").Value = 'Reports tests (Cleaning Jobs) I Me Verify "Download Report" button is disabled by default and the empty PDF message in the container is shown Input worksite, robot and month, then select first report shown, assert API call, download the report then verify the report has been downloaded (Excluding CC) I Me Input company, worksite, robot and month but do not select any report, verify "Download Report" button is disabled and correct empty PDF frame (Excluding CC) Input company, worksite, robot and month, download the first report, verify successful api call and file download user 32131 use213123r 5434634'
$ws.Range("D15").Value = '[''Reports'', ''tests'', ''('', ''Cleaning'', ''Jobs'', '')'', ''I'', ''Me'', ''Verify'', ''``'', ''Download'', ''Report'', "''''", ''button'', ''is'', ''disabled'', ''by'', ''default'', ''and'', ''the'', ''empty'', ''PDF'', ''message'', ''in'', ''the'', ''container'', ''is'', ''shown'', ''Input'', ''worksite'', '','', ''robot'', ''and'', ''month'', '','', ''then'', ''select'', ''first'', ''report'', ''shown'', '','', ''assert'', ''API'', ''call'', '','', ''download'', ''the'', ''report'', ''then'', ''verify'', ''the'', ''report'', ''has'', ''been'', ''downloaded'', ''('', ''Excluding'', ''CC'', '')'', ''I'', ''Me'', ''Input'', ''company'', '','', ''worksite'', '','', ''robot'', ''and'', ''month'', ''but'', ''do'', ''not'', ''select'', ''any'', ''report'', '','', ''verify'', ''``'', ''Download'', ''Report'', "''''", ''button'', ''is'', ''disabled'', ''and'', ''correct'', ''empty'', ''PDF'', ''frame'', ''('', ''Excluding'', ''CC'', '')'', ''Input'', ''company'', '','', ''worksite'', '','', ''robot'', ''and'', ''month'', '','', ''download'', ''the'', ''first'', ''report'', '','', ''verify'', ''successful'', ''api'', ''call'', ''and'', ''file'', ''download'', ''user'', ''32131'', ''use213123r'', ''5434634'']'
